$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 22.4016796005221
$ws.Cells.Item(2, 3).Value = 7.243596312492197
$ws.Cells.Item(2, 4).Value = 13.4931077226102
$ws.Cells.Item(2, 5).Value = 13.72984108201061
$ws.Cells.Item(2, 7).Value = 3.706138319499908
$ws.Cells.Item(2, 10).Value = 8.34187757135391
$ws.Cells.Item(2, 12).Value = 12.78354738944779
$ws.Cells.Item(2, 14).Value = 19.55786906436056
$ws.Cells.Item(2, 15).Value = 32.14037786961303

$ws.Cells.Item(3, 2).Value = 22.0593052593947
$ws.Cells.Item(3, 3).Value = 6.941061014812427
$ws.Cells.Item(3, 4).Value = 13.50005704823367
$ws.Cells.Item(3, 5).Value = 13.76302838162216
$ws.Cells.Item(3, 7).Value = 3.708886957952482
$ws.Cells.Item(3, 10).Value = 8.350449333397798
$ws.Cells.Item(3, 12).Value = 12.77602942142691
$ws.Cells.Item(3, 14).Value = 19.62602733005188
$ws.Cells.Item(3, 15).Value = 32.19849744630404

$ws.Cells.Item(4, 2).Value = 21.85081713257965
$ws.Cells.Item(4, 3).Value = 6.747068331088045
$ws.Cells.Item(4, 4).Value = 13.5067764012555
$ws.Cells.Item(4, 5).Value = 13.78512670062921
$ws.Cells.Item(4, 7).Value = 3.710664429424212
$ws.Cells.Item(4, 10).Value = 8.35599426596775
$ws.Cells.Item(4, 12).Value = 12.77331011943408
$ws.Cells.Item(4, 14).Value = 19.66984256005751
$ws.Cells.Item(4, 15).Value = 32.24151114122029

$ws.Cells.Item(5, 2).Value = 21.76639722047495
$ws.Cells.Item(5, 3).Value = 6.66601167948434
$ws.Cells.Item(5, 4).Value = 13.51013150214801
$ws.Cells.Item(5, 5).Value = 13.79456503422889
$ws.Cells.Item(5, 7).Value = 3.711411419619971
$ws.Cells.Item(5, 10).Value = 8.358324953664841
$ws.Cells.Item(5, 12).Value = 12.77268062678977
$ws.Cells.Item(5, 14).Value = 19.68819352279794
$ws.Cells.Item(5, 15).Value = 32.26087756052208

$ws.Cells.Item(6, 2).Value = 21.75241509027743
$ws.Cells.Item(6, 3).Value = 6.652433474692863
$ws.Cells.Item(6, 4).Value = 13.51072587631202
$ws.Cells.Item(6, 5).Value = 13.79615842815141
$ws.Cells.Item(6, 7).Value = 3.711536827387333
$ws.Cells.Item(6, 10).Value = 8.358716262056367
$ws.Cells.Item(6, 12).Value = 12.77260505426909
$ws.Cells.Item(6, 14).Value = 19.69127068417665
$ws.Cells.Item(6, 15).Value = 32.26420420258565

$ws.Cells.Item(7, 2).Value = 21.84967628783156
$ws.Cells.Item(7, 3).Value = 6.745983186455111
$ws.Cells.Item(7, 4).Value = 13.50681915133334
$ws.Cells.Item(7, 5).Value = 13.78525223535926
$ws.Cells.Item(7, 7).Value = 3.710674411760848
$ws.Cells.Item(7, 10).Value = 8.356025410343729
$ws.Cells.Item(7, 12).Value = 12.77329968978574
$ws.Cells.Item(7, 14).Value = 19.67008803758406
$ws.Cells.Item(7, 15).Value = 32.24176488822212

$ws.Cells.Item(8, 2).Value = 22.28333287038052
$ws.Cells.Item(8, 3).Value = 7.141028893964714
$ws.Cells.Item(8, 4).Value = 13.49499504052048
$ws.Cells.Item(8, 5).Value = 13.7409270032721
$ws.Cells.Item(8, 7).Value = 3.707067458693638
$ws.Cells.Item(8, 10).Value = 8.344774748129046
$ws.Cells.Item(8, 12).Value = 12.78056259284178
$ws.Cells.Item(8, 14).Value = 19.58096299644404
$ws.Cells.Item(8, 15).Value = 32.15889387606231

$ws.Cells.Item(9, 2).Value = 23.14285218621498
$ws.Cells.Item(9, 3).Value = 7.847870596610298
$ws.Cells.Item(9, 4).Value = 13.49124995929332
$ws.Cells.Item(9, 5).Value = 13.66764960017446
$ws.Cells.Item(9, 7).Value = 3.700703259235258
$ws.Cells.Item(9, 10).Value = 8.324938754805157
$ws.Cells.Item(9, 12).Value = 12.80977371612239
$ws.Cells.Item(9, 14).Value = 19.42171332788092
$ws.Cells.Item(9, 15).Value = 32.05472057118972

$ws.Cells.Item(10, 2).Value = 23.77379991157393
$ws.Cells.Item(10, 3).Value = 8.323057394133993
$ws.Cells.Item(10, 4).Value = 13.50031716359604
$ws.Cells.Item(10, 5).Value = 13.62211351342729
$ws.Cells.Item(10, 7).Value = 3.696454871758309
$ws.Cells.Item(10, 10).Value = 8.311708967598454
$ws.Cells.Item(10, 12).Value = 12.84024079033526
$ws.Cells.Item(10, 14).Value = 19.31407287965522
$ws.Cells.Item(10, 15).Value = 32.01398215286687

$ws.Cells.Item(11, 2).Value = 24.05946495125689
$ws.Cells.Item(11, 3).Value = 8.529206366865647
$ws.Cells.Item(11, 4).Value = 13.50699659830215
$ws.Cells.Item(11, 5).Value = 13.60319678031098
$ws.Cells.Item(11, 7).Value = 3.694613941366665
$ws.Cells.Item(11, 10).Value = 8.305979270320289
$ws.Cells.Item(11, 12).Value = 12.85602674644989
$ws.Cells.Item(11, 14).Value = 19.26711488261891
$ws.Cells.Item(11, 15).Value = 32.00326116683395

$ws.Cells.Item(12, 2).Value = 24.16734342669723
$ws.Cells.Item(12, 3).Value = 8.605799996308882
$ws.Cells.Item(12, 4).Value = 13.50989184846302
$ws.Cells.Item(12, 5).Value = 13.59629175840988
$ws.Cells.Item(12, 7).Value = 3.693929934054967
$ws.Cells.Item(12, 10).Value = 8.303850864185023
$ws.Cells.Item(12, 12).Value = 12.8622785611696
$ws.Cells.Item(12, 14).Value = 19.24962024803273
$ws.Cells.Item(12, 15).Value = 32.00032675715562

$ws.Cells.Item(13, 2).Value = 24.14412463939608
$ws.Cells.Item(13, 3).Value = 8.589370017358704
$ws.Cells.Item(13, 4).Value = 13.50925205738668
$ws.Cells.Item(13, 5).Value = 13.59776739209603
$ws.Cells.Item(13, 7).Value = 3.694076665182126
$ws.Cells.Item(13, 10).Value = 8.304307420517079
$ws.Cells.Item(13, 12).Value = 12.86091998395936
$ws.Cells.Item(13, 14).Value = 19.25337527376631
$ws.Cells.Item(13, 15).Value = 32.00090865477574

$ws.Cells.Item(14, 2).Value = 24.06834669287696
$ws.Cells.Item(14, 3).Value = 8.535537397417103
$ws.Cells.Item(14, 4).Value = 13.50722746853228
$ws.Cells.Item(14, 5).Value = 13.60262352373418
$ws.Cells.Item(14, 7).Value = 3.694557405271231
$ws.Cells.Item(14, 10).Value = 8.30580333834313
$ws.Cells.Item(14, 12).Value = 12.85653561320791
$ws.Cells.Item(14, 14).Value = 19.26566983794747
$ws.Cells.Item(14, 15).Value = 32.00299718643819

$ws.Cells.Item(15, 2).Value = 24.02188892455849
$ws.Cells.Item(15, 3).Value = 8.502371048152044
$ws.Cells.Item(15, 4).Value = 13.50603495162963
$ws.Cells.Item(15, 5).Value = 13.60563168146004
$ws.Cells.Item(15, 7).Value = 3.694853578075548
$ws.Cells.Item(15, 10).Value = 8.306725004707632
$ws.Cells.Item(15, 12).Value = 12.85388565209702
$ws.Cells.Item(15, 14).Value = 19.27323799131826
$ws.Cells.Item(15, 15).Value = 32.00442308551685

$ws.Cells.Item(16, 2).Value = 23.75509496395688
$ws.Cells.Item(16, 3).Value = 8.309380826157598
$ws.Cells.Item(16, 4).Value = 13.499931934035
$ws.Cells.Item(16, 5).Value = 13.62338592733722
$ws.Cells.Item(16, 7).Value = 3.696577019714673
$ws.Cells.Item(16, 10).Value = 8.31208920749658
$ws.Cells.Item(16, 12).Value = 12.83924765317244
$ws.Cells.Item(16, 14).Value = 19.31718198198491
$ws.Cells.Item(16, 15).Value = 32.0148401737292

$ws.Cells.Item(17, 2).Value = 23.59100465387992
$ws.Cells.Item(17, 3).Value = 8.188399094462556
$ws.Cells.Item(17, 4).Value = 13.49684128524985
$ws.Cells.Item(17, 5).Value = 13.63473791940476
$ws.Cells.Item(17, 7).Value = 3.697657726907443
$ws.Cells.Item(17, 10).Value = 8.315453750583178
$ws.Cells.Item(17, 12).Value = 12.83075903718967
$ws.Cells.Item(17, 14).Value = 19.34465351509004
$ws.Cells.Item(17, 15).Value = 32.02323307585409

$ws.Cells.Item(18, 2).Value = 23.49650245114759
$ws.Cells.Item(18, 3).Value = 8.117872315361835
$ws.Cells.Item(18, 4).Value = 13.49530419158162
$ws.Cells.Item(18, 5).Value = 13.64143652014715
$ws.Cells.Item(18, 7).Value = 3.698287954876877
$ws.Cells.Item(18, 10).Value = 8.317416122626927
$ws.Cells.Item(18, 12).Value = 12.82605811783958
$ws.Cells.Item(18, 14).Value = 19.36064351867648
$ws.Cells.Item(18, 15).Value = 32.02879558095437

$ws.Cells.Item(19, 2).Value = 23.46448786831713
$ws.Cells.Item(19, 3).Value = 8.093832499155937
$ws.Cells.Item(19, 4).Value = 13.49482511632285
$ws.Cells.Item(19, 5).Value = 13.64373362346591
$ws.Cells.Item(19, 7).Value = 3.698502824261472
$ws.Cells.Item(19, 10).Value = 8.318085221509266
$ws.Cells.Item(19, 12).Value = 12.8244977290127
$ws.Cells.Item(19, 14).Value = 19.36608998199213
$ws.Cells.Item(19, 15).Value = 32.03080513268261

$ws.Cells.Item(20, 2).Value = 23.60848570115998
$ws.Cells.Item(20, 3).Value = 8.201375461262172
$ws.Cells.Item(20, 4).Value = 13.49714540232676
$ws.Cells.Item(20, 5).Value = 13.63351196667229
$ws.Cells.Item(20, 7).Value = 3.6975417907323
$ws.Cells.Item(20, 10).Value = 8.315092778041805
$ws.Cells.Item(20, 12).Value = 12.83164390041131
$ws.Cells.Item(20, 14).Value = 19.34170956231959
$ws.Cells.Item(20, 15).Value = 32.02226353618155

$ws.Cells.Item(21, 2).Value = 24.09061335555335
$ws.Cells.Item(21, 3).Value = 8.551389468113737
$ws.Cells.Item(21, 4).Value = 13.50781222206471
$ws.Cells.Item(21, 5).Value = 13.60119015061822
$ws.Cells.Item(21, 7).Value = 3.694415844930286
$ws.Cells.Item(21, 10).Value = 8.305362831527663
$ws.Cells.Item(21, 12).Value = 12.85781599735071
$ws.Cells.Item(21, 14).Value = 19.26205083945707
$ws.Cells.Item(21, 15).Value = 32.00235317703837

$ws.Cells.Item(22, 2).Value = 24.4039380173567
$ws.Cells.Item(22, 3).Value = 8.77156198140692
$ws.Cells.Item(22, 4).Value = 13.51691545536886
$ws.Cells.Item(22, 5).Value = 13.58157165198495
$ws.Cells.Item(22, 7).Value = 3.692449259196725
$ws.Cells.Item(22, 10).Value = 8.299244422160145
$ws.Cells.Item(22, 12).Value = 12.87651661180814
$ws.Cells.Item(22, 14).Value = 19.2116634406898
$ws.Cells.Item(22, 15).Value = 31.99590082557876

$ws.Cells.Item(23, 2).Value = 24.23690561537838
$ws.Cells.Item(23, 3).Value = 8.654845725726288
$ws.Cells.Item(23, 4).Value = 13.5118623735688
$ws.Cells.Item(23, 5).Value = 13.59190471020624
$ws.Cells.Item(23, 7).Value = 3.693491895605282
$ws.Cells.Item(23, 10).Value = 8.302487974593358
$ws.Cells.Item(23, 12).Value = 12.86639078377783
$ws.Cells.Item(23, 14).Value = 19.23840342590688
$ws.Cells.Item(23, 15).Value = 31.99874375520782

$ws.Cells.Item(24, 2).Value = 23.60058302421486
$ws.Cells.Item(24, 3).Value = 8.195511873524492
$ws.Cells.Item(24, 4).Value = 13.49700716406346
$ws.Cells.Item(24, 5).Value = 13.63406568353282
$ws.Cells.Item(24, 7).Value = 3.697594177701262
$ws.Cells.Item(24, 10).Value = 8.315255886333944
$ws.Cells.Item(24, 12).Value = 12.83124329502574
$ws.Cells.Item(24, 14).Value = 19.34303991189809
$ws.Cells.Item(24, 15).Value = 32.02269956845958

$ws.Cells.Item(25, 2).Value = 22.9099901287771
$ws.Cells.Item(25, 3).Value = 7.664229161329726
$ws.Cells.Item(25, 4).Value = 13.49018425467617
$ws.Cells.Item(25, 5).Value = 13.68601412398398
$ws.Cells.Item(25, 7).Value = 3.702349541074021
$ws.Cells.Item(25, 10).Value = 8.330067966806654
$ws.Cells.Item(25, 12).Value = 12.80028032131476
$ws.Cells.Item(25, 14).Value = 19.46314338332623
$ws.Cells.Item(25, 15).Value = 32.07663049420722
